$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3196676666666667
$ws.Range("H2").Value = 0.959003
$ws.Range("I2").Value = 0.02778181741070332
$ws.Range("J2").Value = 0.02778181741070332
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.6034704469962781
$ws.Range("Q2").Value = 5.893495289764111
$ws.Range("R2").Value = 53.04145760787701
$ws.Range("S2").Value = 0.01676550577120612
$ws.Range("T2").Value = 0.01676550577120611
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3196676666666667
$ws.Range("H3").Value = 0.959003
$ws.Range("I3").Value = 0.02778181741070332
$ws.Range("J3").Value = 0.02778181741070332
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 1.709460674840556
$ws.Range("R3").Value = 15.385146073565
$ws.Range("S3").Value = 0.004862983917110479
$ws.Range("T3").Value = 0.004862983917110479
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3196676666666667
$ws.Range("H4").Value = 0.959003
$ws.Range("I4").Value = 0.02778181741070332
$ws.Range("J4").Value = 0.02778181741070332
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 2.163048848221667
$ws.Range("R4").Value = 19.467439633995
$ws.Range("S4").Value = 0.006153327722386725
$ws.Range("T4").Value = 0.006153327722386725
$ws.Range("G5").Value = 5.787456
$ws.Range("I5").Value = 0.5029787577238425
$ws.Range("J5").Value = 0.5029787577238426
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.6034704469962781
$ws.Range("Q5").Value = 106.699388872768
$ws.Range("R5").Value = 960.294499854912
$ws.Range("S5").Value = 0.3035328157532399
$ws.Range("T5").Value = 0.3035328157532399
$ws.Range("G6").Value = 5.787456
$ws.Range("I6").Value = 0.5029787577238425
$ws.Range("J6").Value = 0.5029787577238426
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("Q6").Value = 30.94910580896
$ws.Range("R6").Value = 278.54195228064
$ws.Range("S6").Value = 0.08804239021875179
$ws.Range("T6").Value = 0.08804239021875179
$ws.Range("G7").Value = 5.787456
$ws.Range("I7").Value = 0.5029787577238425
$ws.Range("J7").Value = 0.5029787577238426
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 39.16113933408
$ws.Range("R7").Value = 352.45025400672
$ws.Range("S7").Value = 0.1114035517518508
$ws.Range("T7").Value = 0.1114035517518508
$ws.Range("G8").Value = 5.399239000000001
$ws.Range("H8").Value = 16.197717
$ws.Range("I8").Value = 0.4692394248654542
$ws.Range("J8").Value = 0.4692394248654542
$ws.Range("M8").Value = 18.43631966666667
$ws.Range("N8").Value = 55.308959
$ws.Range("O8").Value = 0.6034704469962782
$ws.Range("P8").Value = 0.6034704469962781
$ws.Range("Q8").Value = 99.54209616073368
$ws.Range("R8").Value = 895.8788654466031
$ws.Range("S8").Value = 0.2831721254718321
$ws.Range("T8").Value = 0.2831721254718321
$ws.Range("G9").Value = 5.399239000000001
$ws.Range("H9").Value = 16.197717
$ws.Range("I9").Value = 0.4692394248654542
$ws.Range("J9").Value = 0.4692394248654542
$ws.Range("O9").Value = 0.1750419652256785
$ws.Range("P9").Value = 0.1750419652256784
$ws.Range("Q9").Value = 28.87306946244833
$ws.Range("R9").Value = 259.857625162035
$ws.Range("S9").Value = 0.08213659108981619
$ws.Range("T9").Value = 0.08213659108981619
$ws.Range("G10").Value = 5.399239000000001
$ws.Range("H10").Value = 16.197717
$ws.Range("I10").Value = 0.4692394248654542
$ws.Range("J10").Value = 0.4692394248654542
$ws.Range("M10").Value = 6.766555
$ws.Range("N10").Value = 20.299665
$ws.Range("O10").Value = 0.2214875877780434
$ws.Range("P10").Value = 0.2214875877780434
$ws.Range("Q10").Value = 36.534247651645
$ws.Range("R10").Value = 328.808228864805
$ws.Range("S10").Value = 0.1039307083038059
$ws.Range("T10").Value = 0.1039307083038059
